# prob26 (cal part) - add a new y0007/y0008/y0009 building-block group
# to Sheet1, rows 197-199, and move the active selection/view down to
# reflect the newly appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the order the original author entered them (preserves the
# same shared-string append order the commit produced).

# Row 198 description (express the (k/n)*const/n series as a definite integral)
$ws.Range("B198").Value = "`$\dfrac{k}{n}`$를 포함한 식과 `$\dfrac{상수}{n}`$ 의 곱 꼴의 급수를 정적분으로 표현합니다."

# Row 197 description (rewrite the sum as a (k/n)-term times a const/n term)
$ws.Range("B197").Value = "`$\displaystyle\sum`$ 안의 식을 `$\dfrac{k}{n}`$를 포함한 식과 `$\dfrac{상수}{n}`$ 의 곱으로 변형합니다."

# Row 199 formula (the f'(x)/f(x)-type definite integral to evaluate)
$ws.Range("C199").Value = "`$\displaystyle\int_{0}^{1} \dfrac{x^{2}+2 x}{x^{3}+3 x^{2}+1} d x`$;"

# Keys for the new building-block rows
$ws.Range("A197").Value = "y0007"
$ws.Range("A198").Value = "y0008"
$ws.Range("A199").Value = "y0009"

# Row 199 description (evaluate that f'(x)/f(x)-type integral)
$ws.Range("B199").Value = "`$\dfrac{f^{\prime}(x)}{f(x)}`$ 꼴에 대한 정적분을 계산합니다."

# Update the worksheet view: scrolled a few rows further down, with the
# new B201 (blank row following the inserted block) selected.
$ws.Activate()
$ws.Range("B201").Select()
$excel.ActiveWindow.ScrollRow = 187
$excel.ActiveWindow.ScrollColumn = 1
